$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A first (new / renamed test-step identifiers), rows 15-18 ---
$ws.Range("A15").Value = "103_TruckInsurance_001_SmokeTest_FillPageVehicleData"
$ws.Range("A16").Value = "103_TruckInsurance_001_SmokeTest_FillPageInsurantData"
$ws.Range("A17").Value = "103_TruckInsurance_001_SmokeTest_FillPageProductData"
$ws.Range("A18").Value = "103_TruckInsurance_001_SmokeTest_FillPageSendQuote"

# --- Action / Modus (B) + the target selector + NOP (H) for rows 15-18 ---
$ws.Range("B15").Value = "<SET>"
$ws.Range("C15").Value = "103_TruckInsurance_001_SmokeTest_FillPage"
$ws.Range("H15").Value = "<NOP>"

$ws.Range("B16").Value = "<SET>"
$ws.Range("D16").Value = "103_TruckInsurance_001_SmokeTest_FillPage"
$ws.Range("H16").Value = "<NOP>"

$ws.Range("B17").Value = "<SET>"
$ws.Range("E17").Value = "103_TruckInsurance_001_SmokeTest_FillPage"
$ws.Range("H17").Value = "<NOP>"

$ws.Range("B18").Value = "<SET>"
$ws.Range("G18").Value = "103_TruckInsurance_001_SmokeTest_FillPage"
$ws.Range("H18").Value = "<NOP>"

# --- Row 19: "Button Next" step ---
$ws.Range("A19").Value = "Button Next from Page VehicleData"
$ws.Range("B19").Value = "<SET>"
$ws.Range("C19").Value = "Button Next"
$ws.Range("H19").Value = "<NOP>"

# --- Rows 20-23: price-option choices ---
$ws.Range("A20").Value = "Choose Silver"
$ws.Range("B20").Value = "<SET>"
$ws.Range("F20").Value = "Choose Silver"
$ws.Range("H20").Value = "<NOP>"

$ws.Range("A21").Value = "Choose Gold"
$ws.Range("B21").Value = "<SET>"
$ws.Range("F21").Value = "Choose Gold"
$ws.Range("H21").Value = "<NOP>"

$ws.Range("A22").Value = "Choose Platinum"
$ws.Range("B22").Value = "<SET>"
$ws.Range("F22").Value = "Choose Platinum"
$ws.Range("H22").Value = "<NOP>"

$ws.Range("A23").Value = "Choose Ultimate"
$ws.Range("B23").Value = "<SET>"
$ws.Range("F23").Value = "Choose Ultimate"
$ws.Range("H23").Value = "<NOP>"

# --- Row 24: send-quote button ---
$ws.Range("A24").Value = "Send Quote - Button Main Page"
$ws.Range("B24").Value = "<SET>"
$ws.Range("G24").Value = "Button Main Page"
$ws.Range("H24").Value = "<NOP>"

# --- Column widths for C:G (new content is wider than the old bestFit values) ---
$ws.Columns.Item(3).ColumnWidth = 44.053385416666664
$ws.Columns.Item(4).ColumnWidth = 44.053385416666664
$ws.Columns.Item(5).ColumnWidth = 49.166666666666664
$ws.Columns.Item(7).ColumnWidth = 49.166666666666664

# --- Move/resize the picture so it still sits below the (now longer) data table ---
$shp = $ws.Shapes.Item(1)
$shp.Top = 378.0
$shp.Left = 0.6
$shp.Width = 1146.6066929133858
$shp.Height = 719.91

# --- Selection moved to the newly-added rows ---
$ws.Range("A19:XFD24").Select() | Out-Null
